$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capitalize "Data" suffix on the three data-file name cells.
$ws.Range("G8").Value = "SEM1_Data"
$ws.Range("G10").Value = "Fatigue_Data"
$ws.Range("G11").Value = "SEM2_Data"

# Move the active selection/view to G12 (also resets the scrolled
# top-left cell back to the default).
$ws.Range("G12").Select()
